$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

# Update "Scale class" entry on row 2 from "Ordinal" to "Numerical"
$ws.Range("L2").Value = "Numerical"

# Update "Units" entries on rows 2 and 3 from "unit1"/"unit2" to "Index"
$ws.Range("M2").Value = "Index"
$ws.Range("M3").Value = "Index"

# Update selected cell to reflect the active cell after edits
$ws.Range("M3").Select()
